# Applies crypto price/volume updates for Mon May  6 03:36:42 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.068.44"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.145.56"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'591.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").Value = "'146.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.136.02"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("D11").Value = "'5.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("D12").Value = "'0.456"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "'37.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "3.665.61"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "63.820.19"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "3.138.26"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "'468.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").Value = "'14.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.67%  "
$ws.Range("D25").Value = "'13.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "'80.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +10.63%  "
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").Value = "'7.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.90%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D34").Value = "'27.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("D35").Value = "0.0₃0860"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").Value = "'3.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").Value = "'462.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.28%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'51.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'9.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.29%  "
$ws.Range("D43").Value = "'0.292"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.08%  "
$ws.Range("D44").Value = "'0.0373"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").Value = "2.890.29"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").Value = "'39.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.61%  "
$ws.Range("D47").Value = "'0.108"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'132.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.21%  "
$ws.Range("D50").Value = "'0.110"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("E51").Value = "  +4.38%  "
